# Informe de Avance - Periodo 3: "cambios en el periodo 3"
#
# Updates the raw Planned Value (PV) inputs on the "Cálculos" sheet for four
# activities (rows 29-32) so that more planned work has accrued by period 3.
# Every other changed cell in the workbook (period totals, SPI/CPI block,
# the EV chart cache, and the narrative "Análisis" paragraph's SPI figure on
# the "Informe" sheet) is a formula / shared-text consequence of these raw
# inputs, so we only need to touch the inputs themselves plus the narrative
# text, then let the engine recalc.

$wb = $excel.ActiveWorkbook

$wsCalc = $wb.Worksheets.Item("Cálculos")
$wsInforme = $wb.Worksheets.Item("Informe")

# --- Raw Planned Value (PV) input edits, rows 29-32 ------------------------

# Row 29 "Desarrollo vista de agregado de preguntas": Período 2 (F) and
# Período 3 (I) PV go from 4 to 8.
$wsCalc.Range("F29").Value = 8
$wsCalc.Range("I29").Value = 8

# Row 30 "Vista de Administradores": Período 6 (R) and Período 7 (U) PV go
# from blank (0) to 8.
$wsCalc.Range("R30").Value = 8
$wsCalc.Range("U30").Value = 8

# Row 31 "Vista generación de exámen": Período 4 (L) and Período 5 (O) PV go
# from blank (0) to 8.
$wsCalc.Range("L31").Value = 8
$wsCalc.Range("O31").Value = 8

# Row 32 "Vista notas exámenes": Período 6 (R) and Período 7 (U) PV go from
# blank (0) to 8.
$wsCalc.Range("R32").Value = 8
$wsCalc.Range("U32").Value = 8

# --- Narrative "Análisis" paragraph update on the "Informe" sheet ----------
# Cell B48 (merged B48:K50) holds a rich-text paragraph whose body (12pt run)
# quotes the period's SPI value. Only the "0,42" -> "0,39" figure changes;
# the run/formatting structure (bold+underline 14pt "Análisis:" label that
# inherits the cell's base font, a 14pt space run, and the 12pt body run)
# must be preserved, so we surgically replace just that substring via the
# Characters collection instead of overwriting the whole cell value.

$analysisCell = $wsInforme.Range("B48")
$oldText = $analysisCell.Value2
$target = "0,42"
$idxZeroBased = $oldText.IndexOf($target)
if ($idxZeroBased -ge 0) {
    $editRange = $analysisCell.Characters($idxZeroBased + 1, $target.Length)
    $editRange.Text = "0,39"

    # Replacing text via Characters() collapses the rich-text run formatting
    # of the whole string to the cell's default font, so re-apply the
    # original run formatting: "Análisis:" (chars 1-9) keeps the inherited
    # default (no explicit run formatting), the following space (char 10)
    # is 14pt, and the remaining body text is 12pt - neither of the latter
    # two runs are bold/underlined.
    $newLen = $analysisCell.Value2.Length

    $spaceRun = $analysisCell.Characters(10, 1)
    $spaceRun.Font.Size = 14
    $spaceRun.Font.Bold = $false
    $spaceRun.Font.Underline = $false

    $bodyRun = $analysisCell.Characters(11, $newLen - 10)
    $bodyRun.Font.Size = 12
    $bodyRun.Font.Bold = $false
    $bodyRun.Font.Underline = $false
}

# --- View-state housekeeping ------------------------------------------------
# Matches the saved cursor/scroll state: user ends up on "Cálculos" with the
# last-touched cell (V32) selected, and leaves "Informe" with the updated
# paragraph's merged range selected.
$wsInforme.Activate()
$wsInforme.Range("B48:K50").Select()

$wsCalc.Activate()
$wsCalc.Range("V32").Select()
